$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "H2 Oh No"
$ws.Range("C7").Value = 1

$ws.Range("B8").Value = "Unescapable Undine"
$ws.Range("D8").Value = 1

$ws.PageSetup.Orientation = 1

$ws.Range("D11").Select()
